$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

$dateCell = $ws.Cells.Item($row, 4)
$dateCell.Value = 44656
$dateCell.NumberFormat = $ws.Cells.Item(36, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112040
$ws.Cells.Item($row, 7).Value = "Cilantro"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 200
$ws.Cells.Item($row, 11).Value = 600
$ws.Cells.Item($row, 12).Value = 650
$ws.Cells.Item($row, 13).Value = 625
$ws.Cells.Item($row, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item($row, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($row, 16).Value = 625
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
